# Update data paths / re-generate scraped data columns workbook.
# (commit: "da update duong dan de tao data")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header / mapping row (row 2): extra scrape selectors ---
$ws.Range("H2").Value = "href"
$ws.Range("I2").Value = "img class"
$ws.Range("J2").Value = "class"
$ws.Range("L2").Value = "style =width 92%"
$ws.Range("M2").Value = "class review"

# --- New sample/data row (row 4): css-path style selector values ---
$ws.Range("A4").Value = "div"
$ws.Range("B4").Value = "div"
$ws.Range("C4").Value = "div"
$ws.Range("D4").Value = "div"
$ws.Range("E4").Value = "div"
$ws.Range("F4").Value = "div"
$ws.Range("G4").Value = "div"
$ws.Range("H4").Value = "div>a href"
$ws.Range("I4").Value = "div>a>div>img>src"
$ws.Range("J4").Value = "div>a>div>p class= price sales>span class regular"
$ws.Range("K4").Value = "div>a>div>p class= price sales>span class= sales tag"
$ws.Range("L4").Value = "div>a>div>p> span style=width"
$ws.Range("M4").Value = "div>a>div>p class review"

# --- Column widths (approximate the author's saved best-fit widths) ---
$ws.Columns.Item(1).ColumnWidth = 17
$ws.Columns.Item(2).ColumnWidth = 9.5
$ws.Columns.Item(3).ColumnWidth = 7.333333333333333
$ws.Columns.Item(5).ColumnWidth = 5.666666666666667
$ws.Columns.Item(6).ColumnWidth = 8.666666666666666
$ws.Columns.Item(7).ColumnWidth = 10.666666666666666
$ws.Columns.Item(8).ColumnWidth = 10.166666666666666
$ws.Columns.Item(9).ColumnWidth = 14.166666666666666
$ws.Columns.Item(10).ColumnWidth = 36.666666666666664
$ws.Columns.Item(11).ColumnWidth = 36.666666666666664
$ws.Columns.Item(12).ColumnWidth = 23.5
$ws.Columns.Item(13).ColumnWidth = 18.666666666666668

# --- View state: zoom + new selected cell ---
$excel.ActiveWindow.Zoom = 120
$ws.Range("D14").Select() | Out-Null
